$wb = $excel.ActiveWorkbook

# --- Rename existing sheet and add new sheet ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "Fall 2024"

$ws2 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws1)
$ws2.Name = "Spring 2025"

# --- Column widths (match Fall 2024 sheet) ---
$ws2.Columns.Item(1).ColumnWidth = 13.3125
$ws2.Columns.Item(3).ColumnWidth = 85.9453125
$ws2.Columns.Item(5).ColumnWidth = 15.3125

# --- Header row ---
$ws2.Range("A1").Value = "Date:"
$ws2.Range("B1").Value = "Hours:"
$ws2.Range("C1").Value = "Task:"
$hdr = $ws2.Range("A1:C1")
$hdr.Font.Bold = $true
$hdr.HorizontalAlignment = -4108  # xlCenter

# --- Data rows ---
$data = @(
    @{ Y=2025; M=1; D=13; H=0.25; T='Complete team registration form for CEAS Expo' },
    @{ Y=2025; M=1; D=15; H=1; T='Discuss methods of disassembling machine code with advisor' },
    @{ Y=2025; M=1; D=15; H=1; T='Set up Virtual Machine as runtime environment for Hackademia' },
    @{ Y=2025; M=1; D=16; H=2; T='Investigated new server implementation' },
    @{ Y=2025; M=1; D=17; H=0.5; T='Set up Virtual Machine as runtime environment for Hackademia' },
    @{ Y=2025; M=1; D=17; H=0.75; T='Experiment with compiling C code and inspecting output patterns' },
    @{ Y=2025; M=1; D=18; H=0.5; T='Investigated new server implementation' },
    @{ Y=2025; M=1; D=19; H=0.25; T='C code compilation process' },
    @{ Y=2025; M=1; D=22; H=0.75; T='Discuss frontend style library options with advisor' },
    @{ Y=2025; M=1; D=23; H=2.5; T='Team meeting to discuss progress and hold-ups, work on C code compilation process' },
    @{ Y=2025; M=1; D=23; H=0.5; T='Discuss methods of disassembling machine code with advisor' },
    @{ Y=2025; M=1; D=24; H=0.75; T='Met with Pratik to discuss new server implementation' },
    @{ Y=2025; M=1; D=30; H=1.5; T='Team meeting to discuss new server structure, class assignments, and plans for further development' },
    @{ Y=2025; M=1; D=30; H=1.25; T='Implement "compile" API route' },
    @{ Y=2025; M=2; D=2; H=0.25; T='Review test plan' },
    @{ Y=2025; M=2; D=5; H=4.25; T='Implement C code compilation process' },
    @{ Y=2025; M=2; D=6; H=1.75; T='Integrated C code compilation with client changes and discussed future development plans' },
    @{ Y=2025; M=2; D=9; H=0.75; T='User Documentation' },
    @{ Y=2025; M=2; D=9; H=0.5; T='Time Tracking' },
    @{ Y=2025; M=2; D=12; H=2.5; T='Update compilation function to handle multiple functions of C code' },
    @{ Y=2025; M=2; D=13; H=1; T='Meet to sync up with team and discuss upcoming work assignments and homework tasks' },
    @{ Y=2025; M=2; D=13; H=2.75; T='Refactor code for Dynamic Callstack Component' },
    @{ Y=2025; M=2; D=14; H=0.5; T='Work on Expo slide deck' },
    @{ Y=2025; M=2; D=17; H=1; T='Reviewed progress and plans to interpret assembly code with advisor' },
    @{ Y=2025; M=2; D=18; H=2.25; T='Added 32-bit compilation and better error handling' },
    @{ Y=2025; M=2; D=19; H=1; T='Fix issues with library includes and filtering for functions when compiling 32-bit code' },
    @{ Y=2025; M=2; D=20; H=1.25; T='Meeting to sync up with the rest of the team and discuss how to merge divergent branches of work' },
    @{ Y=2025; M=2; D=20; H=0.5; T='Troubleshot and fixed filtering of 32-bit functions' },
    @{ Y=2025; M=2; D=24; H=0.5; T='Add basic information to expo poster' },
    @{ Y=2025; M=2; D=24; H=1; T='Researched assembly instructions and discussed project progress with team' },
    @{ Y=2025; M=2; D=24; H=0.75; T='Discuss methods of emulating assembly instructions and displaying memory values with advisor' },
    @{ Y=2025; M=2; D=26; H=0.25; T='Researched assembly instructions' },
    @{ Y=2025; M=2; D=27; H=1.5; T='Team meeting to sync up on tasks and troubleshoot environments' },
    @{ Y=2025; M=2; D=27; H=1; T='Work on parsing assembly instructions to determine proper implementation' },
    @{ Y=2025; M=3; D=2; H=0.5; T='Reviewed expo poster design' },
    @{ Y=2025; M=3; D=5; H=0.75; T='Experiment with assembly instruction implementation' },
    @{ Y=2025; M=3; D=6; H=1; T='Team meeting to sync up on tasks and determine work to be completed before expo' },
    @{ Y=2025; M=3; D=9; H=3.75; T='Investigated available information for each assembly instruction from the Iced API' },
    @{ Y=2025; M=3; D=10; H=1; T='Met with advisor to discuss best ways of getting operands of each assembly instruction' },
    @{ Y=2025; M=3; D=12; H=2.75; T='Implement case structure to parse and emulate common assembly instructions' },
    @{ Y=2025; M=3; D=13; H=1; T='Meet to sync up with team on tasks and plan work to be done over spring break' },
    @{ Y=2025; M=3; D=18; H=1.25; T='Implemented structures to handle different data type numbers, and unsupported instructions' },
    @{ Y=2025; M=3; D=18; H=0.5; T='Implement POP and ADD instructions' },
    @{ Y=2025; M=3; D=19; H=2.25; T='Investigate Iced handling of memory addresses, and set up instruction pointer' },
    @{ Y=2025; M=3; D=20; H=2.25; T='Refine disassembly process to associate address and state with each instruction' },
    @{ Y=2025; M=3; D=26; H=1.75; T='Met to sync up with team on tasks that had been accomplished over spring break, and work to be done' },
    @{ Y=2025; M=3; D=30; H=1.25; T='Resolve merge conflict' },
    @{ Y=2025; M=3; D=31; H=2.25; T='Implement 32-bit register operations for 64-bit assembly code' },
    @{ Y=2025; M=3; D=31; H=0.75; T='Troubleshoot state updating after each instruction' },
    @{ Y=2025; M=3; D=31; H=1.25; T='Met with advisor to discuss project status and plans before expo' },
    @{ Y=2025; M=4; D=1; H=0.5; T='Resolve merge conflict' },
    @{ Y=2025; M=4; D=3; H=3; T='Unified work from various team members tasks, discussed final action items before the expo' },
    @{ Y=2025; M=4; D=3; H=0.25; T='Time Tracking' },
    @{ Y=2025; M=4; D=6; H=2.5; T='Divided final report work and continued work on associating addresses with each instruction' },
    @{ Y=2025; M=4; D=7; H=3.25; T='Discussed final report work and troubleshot final errors' },
    @{ Y=2025; M=4; D=7; H=3; T='Implemented proper stack visualization for values that are stored relative to the base pointer' },
    @{ Y=2025; M=4; D=14; H=1.5; T='Final self assessment and report updates' }
)

$r = 2
foreach ($row in $data) {
    $d = Get-Date -Year $row.Y -Month $row.M -Day $row.D
    $ws2.Cells.Item($r, 1).Value = $d.Date
    $ws2.Cells.Item($r, 2).Value = $row.H
    $ws2.Cells.Item($r, 3).Value = $row.T
    $r += 1
}
$lastRow = $r - 1

# --- Running total label + formula (row 4, matches Fall 2024 layout) ---
$ws2.Range("E4").Value = "Running total:"
$ws2.Range("F4").Formula = "=SUM(B2:B1000)"

# --- Apply cell styles using a single formatted cell + PasteSpecial(formats) to avoid
#     generating a new style per cell ---

# Style for date column (A2:A<lastRow>): bold, centered, short-date format (matches Fall 2024 style s=5)
$a2 = $ws2.Range("A2")
$a2.NumberFormat = "mm-dd-yy"
$a2.HorizontalAlignment = -4108
$a2.Font.Bold = $true
$a2.Copy() | Out-Null
$ws2.Range("A3:A$lastRow").PasteSpecial(-4122) | Out-Null

# Style for hours/task columns (B/C): centered, default font (matches Fall 2024 style s=1)
$b2 = $ws2.Range("B2")
$b2.HorizontalAlignment = -4108
$b2.Copy() | Out-Null
$ws2.Range("B2:C$lastRow").PasteSpecial(-4122) | Out-Null

# Style for E4 ("Running total:"): bold, left aligned (matches Fall 2024 style s=4)
$e4 = $ws2.Range("E4")
$e4.HorizontalAlignment = -4131  # xlLeft
$e4.Font.Bold = $true

# Style for F4 (sum formula): left aligned, default font (matches Fall 2024 style s=3)
$f4 = $ws2.Range("F4")
$f4.HorizontalAlignment = -4131  # xlLeft

$excel.CalculateFull() | Out-Null

# --- Restore selection/view state on Fall 2024 sheet to match target ---
$ws1.Activate()
$ws1.Range("C21").Select() | Out-Null

Write-Host "Workbook now has $($wb.Worksheets.Count) sheets"
